# Auto-generated Excel COM-interop script to apply numeric cell updates
# across multiple worksheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR)
# as described by the source diff (currentAveragePrice / profit recalculation).

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 747.8
$ws.Range("I2").Value = 580
$ws.Range("K2").Value = 580
$ws.Range("M2").Value = -467
$ws.Range("H17").Value = 1357.08
$ws.Range("J17").Value = 1366.826
$ws.Range("L17").Value = 4100.478
$ws.Range("N17").Value = -4436.478
$ws.Range("H74").Value = 15399.267
$ws.Range("I74").Value = 16230
$ws.Range("K74").Value = 16230
$ws.Range("M74").Value = -15294
$ws.Range("H77").Value = 15399.267
$ws.Range("I77").Value = 16230
$ws.Range("K77").Value = 81150
$ws.Range("M77").Value = -76470
$ws.Range("H111").Value = 4661
$ws.Range("I111").Value = 4661
$ws.Range("K111").Value = 13983
$ws.Range("M111").Value = -10916
$ws.Range("H125").Value = 999.625
$ws.Range("J125").Value = 999.6667
$ws.Range("L125").Value = 8997.0003
$ws.Range("N125").Value = -13917.0003
$ws.Range("H131").Value = 669642.4
$ws.Range("I131").Value = 835936.5
$ws.Range("K131").Value = 2507809.5
$ws.Range("M131").Value = -2502769.5
$ws.Range("H138").Value = 2868.0588
$ws.Range("J138").Value = 3177.5952
$ws.Range("L138").Value = 9532.785600000001
$ws.Range("N138").Value = -19812.7856

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H88").Value = 3798.9
$ws.Range("I88").Value = 2215.6
$ws.Range("J88").Value = 5382.2
$ws.Range("K88").Value = 2215.6
$ws.Range("L88").Value = 5382.2
$ws.Range("M88").Value = -1809.6
$ws.Range("N88").Value = -6194.2
$ws.Range("H91").Value = 3798.9
$ws.Range("I91").Value = 2215.6
$ws.Range("J91").Value = 5382.2
$ws.Range("K91").Value = 2215.6
$ws.Range("L91").Value = 5382.2
$ws.Range("M91").Value = -811.5999999999999
$ws.Range("N91").Value = -8190.2
$ws.Range("H122").Value = 2684.4583
$ws.Range("I122").Value = 2782.2856
$ws.Range("K122").Value = 8346.856800000001
$ws.Range("M122").Value = -5896.856800000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H82").Value = 28187.834
$ws.Range("I82").Value = 3165.5
$ws.Range("J82").Value = 78232.5
$ws.Range("K82").Value = 3165.5
$ws.Range("L82").Value = 78232.5
$ws.Range("M82").Value = -2782.5
$ws.Range("N82").Value = -78998.5
$ws.Range("H85").Value = 28187.834
$ws.Range("I85").Value = 3165.5
$ws.Range("J85").Value = 78232.5
$ws.Range("K85").Value = 3165.5
$ws.Range("L85").Value = 78232.5
$ws.Range("M85").Value = -1839.5
$ws.Range("N85").Value = -80884.5
$ws.Range("H105").Value = 12382784
$ws.Range("I105").Value = 770905.7
$ws.Range("K105").Value = 770905.7
$ws.Range("M105").Value = -769158.7

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H2").Value = 82
$ws.Range("I2").Value = 82
$ws.Range("K2").Value = 82
$ws.Range("M2").Value = 31
$ws.Range("H5").Value = 1043
$ws.Range("I5").Value = 1043
$ws.Range("K5").Value = 1043
$ws.Range("M5").Value = -931
$ws.Range("H10").Value = 3543.75
$ws.Range("I10").Value = 1391.6666
$ws.Range("J10").Value = 10000
$ws.Range("K10").Value = 1391.6666
$ws.Range("L10").Value = 10000
$ws.Range("M10").Value = -1252.6666
$ws.Range("N10").Value = -10278
$ws.Range("H12").Value = 2550
$ws.Range("I12").Value = 550
$ws.Range("J12").Value = 4550
$ws.Range("K12").Value = 550
$ws.Range("L12").Value = 4550
$ws.Range("M12").Value = -380
$ws.Range("N12").Value = -4890
$ws.Range("H13").Value = 6247.5
$ws.Range("J13").Value = 6247.5
$ws.Range("L13").Value = 6247.5
$ws.Range("N13").Value = -6525.5
$ws.Range("H15").Value = 15000
$ws.Range("J15").Value = 15000
$ws.Range("L15").Value = 15000
$ws.Range("N15").Value = -15340
$ws.Range("H17").Value = 1999.5
$ws.Range("I17").Value = 1999.5
$ws.Range("K17").Value = 1999.5
$ws.Range("M17").Value = -1825.5
$ws.Range("H132").Value = 19614584
$ws.Range("I132").Value = 4560.5
$ws.Range("K132").Value = 13681.5
$ws.Range("M132").Value = -11151.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H14").Value = 321
$ws.Range("I14").Value = 321
$ws.Range("K14").Value = 963
$ws.Range("M14").Value = -790
$ws.Range("H39").Value = 8448.909
$ws.Range("I39").Value = 1000
$ws.Range("J39").Value = 11242.25
$ws.Range("K39").Value = 3000
$ws.Range("L39").Value = 33726.75
$ws.Range("M39").Value = -2706
$ws.Range("N39").Value = -34314.75
$ws.Range("H82").Value = 8875
$ws.Range("J82").Value = 11500
$ws.Range("L82").Value = 34500
$ws.Range("N82").Value = -35312
$ws.Range("H85").Value = 8875
$ws.Range("J85").Value = 11500
$ws.Range("L85").Value = 34500
$ws.Range("N85").Value = -37308
$ws.Range("H87").Value = 760
$ws.Range("I87").Value = 760
$ws.Range("K87").Value = 2280
$ws.Range("M87").Value = -1032
$ws.Range("H90").Value = 760
$ws.Range("I90").Value = 760
$ws.Range("K90").Value = 6840
$ws.Range("M90").Value = -600
$ws.Range("H94").Value = 4966.6665
$ws.Range("I94").Value = 2450
$ws.Range("K94").Value = 7350
$ws.Range("M94").Value = -6674

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H11").Value = 10552165
$ws.Range("I11").Value = 42500576
$ws.Range("J11").Value = 2565062.5
$ws.Range("K11").Value = 42500576
$ws.Range("L11").Value = 2565062.5
$ws.Range("M11").Value = -42500437
$ws.Range("N11").Value = -2565340.5
$ws.Range("H21").Value = 24001.2
$ws.Range("J21").Value = 24001.2
$ws.Range("L21").Value = 24001.2
$ws.Range("N21").Value = -24347.2
$ws.Range("H30").Value = 24001.2
$ws.Range("J30").Value = 24001.2
$ws.Range("L30").Value = 24001.2
$ws.Range("N30").Value = -24211.2
$ws.Range("H80").Value = 52635148
$ws.Range("J80").Value = 4341.4287
$ws.Range("L80").Value = 4341.4287
$ws.Range("N80").Value = -6337.4287
$ws.Range("H83").Value = 52635148
$ws.Range("J83").Value = 4341.4287
$ws.Range("L83").Value = 21707.1435
$ws.Range("N83").Value = -31691.1435
$ws.Range("H138").Value = 106992
$ws.Range("J138").Value = 106992
$ws.Range("L138").Value = 106992
$ws.Range("N138").Value = -117272

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 1560.1578
$ws.Range("I61").Value = 1527.375
$ws.Range("J61").Value = 1735
$ws.Range("K61").Value = 1527.375
$ws.Range("L61").Value = 1735
$ws.Range("M61").Value = -1325.375
$ws.Range("N61").Value = -2139
$ws.Range("H82").Value = 897.2778
$ws.Range("I82").Value = 896.53845
$ws.Range("J82").Value = 899.2
$ws.Range("K82").Value = 896.53845
$ws.Range("L82").Value = 899.2
$ws.Range("M82").Value = -535.53845
$ws.Range("N82").Value = -1621.2
$ws.Range("H85").Value = 897.2778
$ws.Range("I85").Value = 896.53845
$ws.Range("J85").Value = 899.2
$ws.Range("K85").Value = 896.53845
$ws.Range("L85").Value = 899.2
$ws.Range("M85").Value = 351.46155
$ws.Range("N85").Value = -3395.2
$ws.Range("H113").Value = 1560.1578
$ws.Range("I113").Value = 1527.375
$ws.Range("J113").Value = 1735
$ws.Range("K113").Value = 1527.375
$ws.Range("L113").Value = 1735
$ws.Range("M113").Value = 642.625
$ws.Range("N113").Value = -6075

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H69").Value = 12250
$ws.Range("J69").Value = 12250
$ws.Range("L69").Value = 12250
$ws.Range("N69").Value = -13748
$ws.Range("H72").Value = 12250
$ws.Range("J72").Value = 12250
$ws.Range("L72").Value = 36750
$ws.Range("N72").Value = -44238
$ws.Range("H86").Value = 49999
$ws.Range("J86").Value = 49999
$ws.Range("L86").Value = 49999
$ws.Range("N86").Value = -52245
$ws.Range("H89").Value = 49999
$ws.Range("J89").Value = 49999
$ws.Range("L89").Value = 249995
$ws.Range("N89").Value = -261227
$ws.Range("H107").Value = 1104
$ws.Range("I107").Value = 892.2273
$ws.Range("J107").Value = 1769.5714
$ws.Range("K107").Value = 2676.6819
$ws.Range("L107").Value = 5308.7142
$ws.Range("M107").Value = -756.6819
$ws.Range("N107").Value = -9148.7142
